$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row D/E value updates ---
$ws.Range("D2").Value = '36.790.33'
$ws.Range("E2").Value = '  -1.07%  '

$ws.Range("D3").Value = '2.092.21'
$ws.Range("E3").Value = '  +1.87%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.31'
$ws.Range("E5").Value = '  -0.97%  '

$ws.Range("E6").Value = '  -2.19%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.70'
$ws.Range("E8").Value = '  -4.73%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.54'
$ws.Range("E9").Value = '  -1.25%  '

$ws.Range("E10").Value = '  -4.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0769'
$ws.Range("E11").Value = '  -2.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.110'
$ws.Range("E12").Value = '  +1.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.99'
$ws.Range("E13").Value = '  -5.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.885'
$ws.Range("E14").Value = '  +6.09%  '

$ws.Range("D15").Value = '2.395.51'
$ws.Range("E15").Value = '  +1.83%  '

$ws.Range("E16").Value = '  -3.98%  '

$ws.Range("D17").Value = '2.090.07'

$ws.Range("D18").Value = '36.784.43'
$ws.Range("E18").Value = '  -0.97%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.45'
$ws.Range("E19").Value = '  -3.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.28'
$ws.Range("E20").Value = '  -2.56%  '

$ws.Range("D21").Value = '0.0₃0880'

$ws.Range("E22").Value = '  +1.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.93'
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("E25").Value = '  -2.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.94'
$ws.Range("E26").Value = '  +5.71%  '

$ws.Range("E27").Value = '  -0.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.36'
$ws.Range("E28").Value = '  -0.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.92'
$ws.Range("E29").Value = '  +3.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.36'
$ws.Range("E30").Value = '  +10.74%  '

$ws.Range("E31").Value = '  -0.58%  '

$ws.Range("E32").Value = '  +6.97%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.73'
$ws.Range("E33").Value = '  +4.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0612'
$ws.Range("E34").Value = '  -1.71%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.43'
$ws.Range("E35").Value = '  +6.70%  '

$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("E37").Value = '  +4.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0842'
$ws.Range("E38").Value = '  -6.31%  '

$ws.Range("E39").Value = '  -3.73%  '

$ws.Range("E42").Value = '  -0.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.94'
$ws.Range("E45").Value = '  +0.78%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.44'
$ws.Range("E46").Value = '  -5.21%  '

$ws.Range("D47").Value = '1.341.35'
$ws.Range("E47").Value = '  +4.59%  '

$ws.Range("E48").Value = '  -1.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.07'
$ws.Range("E49").Value = '  +3.16%  '

$ws.Range("E50").Value = '  -0.95%  '

$ws.Range("D51").Value = '2.278.72'
$ws.Range("E51").Value = '  +1.67%  '

# --- Row 40 & 41: content swap (THORChain <-> ARBITRUM) with updated values ---
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.17'
$ws.Range("E40").Value = '  +2.09%  '

$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.93'
$ws.Range("E41").Value = '  -4.91%  '

# --- Row 43 & 44: content swap (HuobiToken <-> Cronos) with updated values ---
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0955'
$ws.Range("E43").Value = '  -8.43%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.89'
$ws.Range("E44").Value = '  -8.04%  '
